$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.862.11"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.637.79"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'215.41"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'0.5050"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "'0.06429"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").Value = "'19.78"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").Value = "'0.07790"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "'4.290"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "1.863.85"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "1.637.31"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "'0.5603"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "0.0₅7626"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "'62.99"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "25.875.25"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'194.54"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'4.329"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D22").Value = "'9.880"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "'6.097"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'1.775"
$ws.Range("E25").Value = "  -6.64%  "
$ws.Range("D26").Value = "'140.35"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "'0.1258"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").Value = "'6.826"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'0.04896"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "'3.295"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").Value = "'3.226"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "'1.568"
$ws.Range("D35").Value = "'2.381"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Value = "'0.9037"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "'2.577"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'0.5515"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "1.127.09"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "'0.01563"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").Value = "'5.550"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'0.8007"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").Value = "'98.17"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "1.774.34"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "0.0₈114"
$ws.Range("E46").Value = "  -5.55%  "
$ws.Range("D47").Value = "'55.38"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "'0.4264"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("D49").Value = "'7.720"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").Value = "'0.05041"
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  +0.44%  "
